# Regenerated save_data: column G ("K" = strikeouts, previously derived from
# a different "Strike#" source) is recalculated for each game row (rows 2-32).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 6
    4  = 3
    5  = 2
    6  = 2
    7  = 3
    8  = 7
    9  = 1
    10 = 7
    11 = 6
    12 = 2
    13 = 5
    14 = 7
    15 = 4
    16 = 4
    17 = 4
    18 = 6
    19 = 4
    20 = 6
    21 = 6
    22 = 4
    23 = 1
    24 = 2
    25 = 2
    26 = 8
    27 = 2
    28 = 6
    29 = 5
    30 = 1
    31 = 1
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
